# Auto-generated Excel COM-interop script
# Combines b10brf and b10brm into b10br: updates summary-statistics rows for
# MHC="Kb" at levels L3/L4/L5/L6 across CDR3/CDR3a/CDR3b/FullChain/FullAlpha/FullBeta
# with recomputed N_Values, Mean/Median/Std/Min/Max, ratios and fold-change labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = [double]"3"
$ws.Range("E4").Value = [double]"0.0001896700925827139"
$ws.Range("F4").Value = [double]"0"
$ws.Range("G4").Value = [double]"0.0003285182370295535"
$ws.Range("I4").Value = [double]"0.0005690102777481419"
$ws.Range("J4").Value = [double]"0.04909701448755735"
$ws.Range("K4").Value = "0.049x"
$ws.Range("L4").Value = [double]"0.001537083567199424"
$ws.Range("M4").Value = "0.002x"
# Row 5
$ws.Range("D5").Value = [double]"6"
$ws.Range("E5").Value = [double]"0.002084510491973922"
$ws.Range("F5").Value = [double]"0.001770259869524957"
$ws.Range("G5").Value = [double]"0.001926021526344725"
$ws.Range("H5").Value = [double]"0.0001263573222030195"
$ws.Range("J5").Value = [double]"10.99019072321526"
$ws.Range("K5").Value = "10.990x"
$ws.Range("L5").Value = [double]"0.01689284156104172"
$ws.Range("M5").Value = "0.017x"
# Row 6
$ws.Range("D6").Value = [double]"15"
$ws.Range("E6").Value = [double]"0.0008717382153061114"
$ws.Range("F6").Value = [double]"0"
$ws.Range("G6").Value = [double]"0.001547896233658058"
$ws.Range("I6").Value = [double]"0.004387668199625277"
$ws.Range("J6").Value = [double]"0.4181980463339483"
$ws.Range("K6").Value = "0.418x"
$ws.Range("L6").Value = [double]"0.007064553337856572"
$ws.Range("M6").Value = "0.007x"
# Row 7
$ws.Range("J7").Value = [double]"1.704861933219134e-05"
# Row 16
$ws.Range("D16").Value = [double]"3"
$ws.Range("E16").Value = [double]"0.001965914872858547"
$ws.Range("F16").Value = [double]"0.0009513140581101746"
$ws.Range("G16").Value = [double]"0.002129071249357088"
$ws.Range("H16").Value = [double]"0.0005339155935228867"
$ws.Range("I16").Value = [double]"0.004412514966942579"
$ws.Range("J16").Value = [double]"0.07273015474656162"
$ws.Range("K16").Value = "0.073x"
$ws.Range("L16").Value = [double]"0.01280550636805721"
$ws.Range("M16").Value = "0.013x"
# Row 17
$ws.Range("D17").Value = [double]"6"
$ws.Range("E17").Value = [double]"0.005714316625944758"
$ws.Range("F17").Value = [double]"0.006886832588299329"
$ws.Range("G17").Value = [double]"0.004041263951423072"
$ws.Range("H17").Value = [double]"0.0003618693134822167"
$ws.Range("I17").Value = [double]"0.009676759149110353"
$ws.Range("J17").Value = [double]"2.906695861980958"
$ws.Range("K17").Value = "2.907x"
$ws.Range("L17").Value = [double]"0.03722171237060268"
$ws.Range("M17").Value = "0.037x"
# Row 18
$ws.Range("D18").Value = [double]"15"
$ws.Range("E18").Value = [double]"0.002986583860895246"
$ws.Range("F18").Value = [double]"0.001203791263138444"
$ws.Range("G18").Value = [double]"0.003485569988878292"
$ws.Range("H18").Value = [double]"8.38635049161365e-05"
$ws.Range("I18").Value = [double]"0.009676759149110353"
$ws.Range("J18").Value = [double]"0.5226493483639382"
$ws.Range("K18").Value = "0.523x"
$ws.Range("L18").Value = [double]"0.01945390371548543"
$ws.Range("M18").Value = "0.019x"
# Row 19
$ws.Range("J19").Value = [double]"0.03348690001099775"
$ws.Range("K19").Value = "0.033x"
# Row 28
$ws.Range("D28").Value = [double]"3"
$ws.Range("E28").Value = [double]"0.006974122396963326"
$ws.Range("F28").Value = [double]"0.001043260971419647"
$ws.Range("G28").Value = [double]"0.0103715879090439"
$ws.Range("H28").Value = [double]"0.0009290870941356379"
$ws.Range("I28").Value = [double]"0.01895001912533469"
$ws.Range("J28").Value = [double]"0.5292522714994019"
$ws.Range("K28").Value = "0.529x"
$ws.Range("L28").Value = [double]"0.04518105962603398"
$ws.Range("M28").Value = "0.045x"
# Row 29
$ws.Range("D29").Value = [double]"6"
$ws.Range("E29").Value = [double]"0.006510989066294701"
$ws.Range("F29").Value = [double]"0.005906672137289037"
$ws.Range("G29").Value = [double]"0.006451943322515748"
$ws.Range("H29").Value = [double]"0.0001725524722557364"
$ws.Range("J29").Value = [double]"0.9335926007162876"
$ws.Range("K29").Value = "0.934x"
$ws.Range("L29").Value = [double]"0.04218070295938672"
$ws.Range("M29").Value = "0.042x"
# Row 30
$ws.Range("D30").Value = [double]"15"
$ws.Range("E30").Value = [double]"0.004697381448048906"
$ws.Range("F30").Value = [double]"0.001272563071040687"
$ws.Range("G30").Value = [double]"0.006209460863732814"
$ws.Range("I30").Value = [double]"0.01895001912533469"
$ws.Range("J30").Value = [double]"0.7214543597324931"
$ws.Range("K30").Value = "0.721x"
$ws.Range("L30").Value = [double]"0.03043145204663082"
$ws.Range("M30").Value = "0.030x"
# Row 31
$ws.Range("J31").Value = [double]"0.002897051496927236"
$ws.Range("K31").Value = "0.003x"
# Row 40
$ws.Range("D40").Value = [double]"3"
$ws.Range("E40").Value = [double]"0.0001896700925827139"
$ws.Range("F40").Value = [double]"0"
$ws.Range("G40").Value = [double]"0.0003285182370295535"
$ws.Range("I40").Value = [double]"0.0005690102777481419"
$ws.Range("J40").Value = [double]"0.06863537683954972"
$ws.Range("K40").Value = "0.069x"
$ws.Range("L40").Value = [double]"0.001594230472176995"
$ws.Range("M40").Value = "0.002x"
# Row 41
$ws.Range("D41").Value = [double]"6"
$ws.Range("E41").Value = [double]"0.001899085284935418"
$ws.Range("F41").Value = [double]"0.001603394922581322"
$ws.Range("G41").Value = [double]"0.001817288391338396"
$ws.Range("H41").Value = [double]"2.717361767806872e-05"
$ws.Range("J41").Value = [double]"10.01257108633107"
$ws.Range("K41").Value = "10.013x"
$ws.Range("L41").Value = [double]"0.01596234593066731"
$ws.Range("M41").Value = "0.016x"
# Row 42
$ws.Range("D42").Value = [double]"15"
$ws.Range("E42").Value = [double]"0.0007975681324907098"
$ws.Range("G42").Value = [double]"0.001437604252366915"
$ws.Range("I42").Value = [double]"0.004237778913302674"
$ws.Range("J42").Value = [double]"0.4199748893940973"
$ws.Range("K42").Value = "0.420x"
$ws.Range("L42").Value = [double]"0.00670378446670232"
$ws.Range("M42").Value = "0.007x"
# Row 43
$ws.Range("J43").Value = [double]"1.531401311139329e-05"
# Row 52
$ws.Range("D52").Value = [double]"3"
$ws.Range("E52").Value = [double]"0.001016998428862044"
$ws.Range("F52").Value = [double]"0.0008090614886731392"
$ws.Range("G52").Value = [double]"0.0008922915950490949"
$ws.Range("H52").Value = [double]"0.0002470355731225296"
$ws.Range("I52").Value = [double]"0.001994898224790463"
$ws.Range("J52").Value = [double]"0.05545135349828925"
$ws.Range("K52").Value = "0.055x"
$ws.Range("L52").Value = [double]"0.007121753942591101"
$ws.Range("M52").Value = "0.007x"
# Row 53
$ws.Range("D53").Value = [double]"6"
$ws.Range("E53").Value = [double]"0.004066272494482025"
$ws.Range("F53").Value = [double]"0.004746230806991261"
$ws.Range("G53").Value = [double]"0.003102102291150781"
$ws.Range("H53").Value = [double]"0.0002326302729528536"
$ws.Range("I53").Value = [double]"0.007119741100323625"
$ws.Range("J53").Value = [double]"3.998307548057792"
$ws.Range("K53").Value = "3.998x"
$ws.Range("L53").Value = [double]"0.02847496254407234"
$ws.Range("M53").Value = "0.028x"
# Row 54
$ws.Range("D54").Value = [double]"15"
$ws.Range("E54").Value = [double]"0.002000508003911323"
$ws.Range("F54").Value = [double]"0.0005698228764310819"
$ws.Range("G54").Value = [double]"0.002591057975854067"
$ws.Range("H54").Value = [double]"4.048582995951417e-05"
$ws.Range("I54").Value = [double]"0.007119741100323625"
$ws.Range("J54").Value = [double]"0.4919758837180819"
$ws.Range("K54").Value = "0.492x"
$ws.Range("L54").Value = [double]"0.01400899486145927"
$ws.Range("M54").Value = "0.014x"
# Row 55
$ws.Range("J55").Value = [double]"0.008304308435942646"
$ws.Range("K55").Value = "0.008x"
# Row 64
$ws.Range("D64").Value = [double]"3"
$ws.Range("E64").Value = [double]"0.006839099823063106"
$ws.Range("F64").Value = [double]"0.0009246417013407304"
$ws.Range("G64").Value = [double]"0.01025142135592126"
$ws.Range("H64").Value = [double]"0.0009162371804883127"
$ws.Range("I64").Value = [double]"0.01867642058736028"
$ws.Range("J64").Value = [double]"0.5304752443202754"
$ws.Range("K64").Value = "0.530x"
$ws.Range("L64").Value = [double]"0.04440528612712769"
$ws.Range("M64").Value = "0.044x"
# Row 65
$ws.Range("D65").Value = [double]"6"
$ws.Range("E65").Value = [double]"0.00631716992090469"
$ws.Range("F65").Value = [double]"0.005880460248546797"
$ws.Range("G65").Value = [double]"0.006238323842731126"
$ws.Range("H65").Value = [double]"0.0001725524722557364"
$ws.Range("J65").Value = [double]"0.9236844152503314"
$ws.Range("K65").Value = "0.924x"
$ws.Range("L65").Value = [double]"0.04101647075035959"
$ws.Range("M65").Value = "0.041x"
# Row 66
$ws.Range("D66").Value = [double]"15"
$ws.Range("E66").Value = [double]"0.004408986963751724"
$ws.Range("F66").Value = [double]"0.001225648303444717"
$ws.Range("G66").Value = [double]"0.006056380445464147"
$ws.Range("I66").Value = [double]"0.01867642058736028"
$ws.Range("J66").Value = [double]"0.6979370539268804"
$ws.Range("K66").Value = "0.698x"
$ws.Range("L66").Value = [double]"0.02862691475798404"
$ws.Range("M66").Value = "0.029x"
# Row 67
$ws.Range("J67").Value = [double]"0.0015704349679927"
$ws.Range("K67").Value = "0.002x"
